$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A3').Value = 'stillen'
$ws.Range('B3').Value = 'flower/flower024.jpg'
$ws.Range('C3').Value = 'flower'
$ws.Range('A4').Value = 'jubeln'
$ws.Range('B4').Value = 'face/face008.jpg'
$ws.Range('C4').Value = 'face'
$ws.Range('A6').Value = 'danken'
$ws.Range('B6').Value = 'flower/flower010.jpg'
$ws.Range('C6').Value = 'flower'
$ws.Range('A7').Value = 'regnen'
$ws.Range('B7').Value = 'flower/flower005.jpg'
$ws.Range('A9').Value = 'pflegen'
$ws.Range('B9').Value = 'face/face016.jpg'
$ws.Range('A10').Value = 'schultern'
$ws.Range('B10').Value = 'flower/flower023.jpg'
$ws.Range('C10').Value = 'flower'
$ws.Range('A12').Value = 'wohnen'
$ws.Range('B12').Value = 'face/face021.jpg'
$ws.Range('A13').Value = 'ehren'
$ws.Range('B13').Value = 'face/face029.jpg'
$ws.Range('A15').Value = 'buchen'
$ws.Range('B15').Value = 'flower/flower026.jpg'
$ws.Range('C15').Value = 'flower'
$ws.Range('A16').Value = 'parken'
$ws.Range('A18').Value = 'erben'
$ws.Range('B18').Value = 'face/face028.jpg'
$ws.Range('C18').Value = 'face'
$ws.Range('A19').Value = 'husten'
$ws.Range('B19').Value = 'face/face001.jpg'
$ws.Range('A21').Value = 'sehen'
$ws.Range('B21').Value = 'flower/flower001.jpg'
$ws.Range('C21').Value = 'flower'
$ws.Range('A22').Value = 'herrschen'
$ws.Range('B22').Value = 'face/face009.jpg'
$ws.Range('A24').Value = 'wandern'
$ws.Range('B24').Value = 'face/face026.jpg'
$ws.Range('C24').Value = 'face'
$ws.Range('A25').Value = 'landen'
$ws.Range('B25').Value = 'flower/flower031.jpg'
$ws.Range('A27').Value = 'heilen'
$ws.Range('B27').Value = 'face/face018.jpg'
$ws.Range('A28').Value = 'tauschen'
$ws.Range('B28').Value = 'face/face017.jpg'
$ws.Range('A30').Value = 'tropfen'
$ws.Range('B30').Value = 'face/face010.jpg'
$ws.Range('C30').Value = 'face'
$ws.Range('A31').Value = 'planen'
$ws.Range('B31').Value = 'flower/flower033.jpg'
$ws.Range('A33').Value = 'holen'
$ws.Range('B33').Value = 'flower/flower003.jpg'
$ws.Range('A34').Value = 'albern'
$ws.Range('B34').Value = 'flower/flower008.jpg'
$ws.Range('A36').Value = 'segeln'
$ws.Range('B36').Value = 'flower/flower000.jpg'
$ws.Range('A37').Value = 'stärken'
$ws.Range('B37').Value = 'face/face013.jpg'
$ws.Range('C37').Value = 'face'
$ws.Range('A39').Value = 'lehnen'
$ws.Range('B39').Value = 'flower/flower022.jpg'
$ws.Range('C39').Value = 'flower'
$ws.Range('A40').Value = 'bauen'
$ws.Range('B40').Value = 'flower/flower004.jpg'
$ws.Range('C40').Value = 'flower'
$ws.Range('A42').Value = 'passen'
$ws.Range('B42').Value = 'flower/flower007.jpg'
$ws.Range('C42').Value = 'flower'
$ws.Range('A43').Value = 'dürfen'
$ws.Range('B43').Value = 'flower/flower032.jpg'
$ws.Range('A45').Value = 'brauchen'
$ws.Range('B45').Value = 'face/face012.jpg'
$ws.Range('C45').Value = 'face'
$ws.Range('A46').Value = 'fühlen'
$ws.Range('B46').Value = 'face/face019.jpg'
$ws.Range('C46').Value = 'face'
$ws.Range('A48').Value = 'füttern'
$ws.Range('B48').Value = 'face/face023.jpg'
$ws.Range('C48').Value = 'face'
$ws.Range('A49').Value = 'werden'
